# Append two new paragraphs at the end of the document body:
#   1) an empty paragraph (en-US language mark only, no run)
#   2) a paragraph containing the text "I want to change this file by develop"
#
# We use Range.InsertXML with a minimal Flat-OPC "pkg:package" wrapper so that
# the freshly-created empty paragraph does not pick up a spurious empty
# <w:r> run the way TypeParagraph()/InsertParagraphAfter() would.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Collapse(0)   # wdCollapseEnd -> collapse to the very end of the document

$flatOpcXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
          '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>' +
          '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
            '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I want to change this file by develop</w:t></w:r>' +
          '</w:p>' +
        '</w:body>' +
      '</w:document>' +
    '</pkg:xmlData>' +
  '</pkg:part>' +
'</pkg:package>'

$rng.InsertXML($flatOpcXml)
